$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "556.24") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '65.261.30'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '3.381.09'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '556.24'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '174.92'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("D8").Value = '3.370.11'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +5.08%  '
$ws.Range("D11").Value = '0.635'
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").Value = '53.77'
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '9.19'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '3.915.03'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '18.34'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D18").Value = '3.355.91'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '65.257.61'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = '11.86'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("D22").Value = '458.48'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("E23").Value = '  +2.58%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '14.10'
$ws.Range("E25").Value = '  +5.62%  '
$ws.Range("D26").Value = '87.69'
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").Value = '2.90'
$ws.Range("E27").Value = '  +2.27%  '
$ws.Range("D28").Value = '10.71'
$ws.Range("E28").Value = '  -1.86%  '
$ws.Range("D29").Value = '8.72'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("D31").Value = '6.55'
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").Value = '63.29'
$ws.Range("E32").Value = '  +7.42%  '
$ws.Range("D33").Value = '11.46'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("D34").Value = '578.03'
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("D39").Value = '35.75'
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("D41").Value = '0.0₃0740'
$ws.Range("E41").Value = '  -2.53%  '
$ws.Range("D42").Value = '3.091.52'
$ws.Range("D43").Value = '0.0418'
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = '2.77'
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("E45").Value = '  -0.61%  '
$ws.Range("E46").Value = '  -2.86%  '
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '140.72'
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("E50").Value = '  -1.86%  '
$ws.Range("D51").Value = '8.31'
$ws.Range("E51").Value = '  -0.75%  '

# Restore the original (default/General) cell formatting now that the
# values are safely stored as text.
$ws.Range("D2:E51").ClearFormats()
